# Fix total marks error on the "quiz" marksheet for roll number 1401CB27.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right/Wrong counts corrected
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): recalculated totals and score summary
$ws.Range("B12").Value = 52
$ws.Range("C12").Value = -10
$ws.Range("E12").Value = "42 / 112"
